$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 379.01  # H15
$ws.Cells.Item(15, 9).Value = 379.01  # I15
$ws.Cells.Item(15, 11).Value = 1137.03  # K15
$ws.Cells.Item(15, 13).Value = -968.03  # M15

$ws.Cells.Item(53, 8).Value = 347.66666  # H53
$ws.Cells.Item(53, 9).Value = 292.75  # I53
$ws.Cells.Item(53, 10).Value = 375.125  # J53
$ws.Cells.Item(53, 11).Value = 292.75  # K53
$ws.Cells.Item(53, 12).Value = 375.125  # L53
$ws.Cells.Item(53, 13).Value = 344.25  # M53
$ws.Cells.Item(53, 14).Value = -1649.125  # N53

$ws.Cells.Item(69, 8).Value = 62503428  # H69
$ws.Cells.Item(69, 9).Value = 2980  # I69
$ws.Cells.Item(69, 10).Value = 66670124  # J69
$ws.Cells.Item(69, 11).Value = 8940  # K69
$ws.Cells.Item(69, 12).Value = 200010372  # L69
$ws.Cells.Item(69, 13).Value = -8066  # M69
$ws.Cells.Item(69, 14).Value = -200012120  # N69

$ws.Cells.Item(72, 8).Value = 62503428  # H72
$ws.Cells.Item(72, 9).Value = 2980  # I72
$ws.Cells.Item(72, 10).Value = 66670124  # J72
$ws.Cells.Item(72, 11).Value = 26820  # K72
$ws.Cells.Item(72, 12).Value = 600031116  # L72
$ws.Cells.Item(72, 13).Value = -22452  # M72
$ws.Cells.Item(72, 14).Value = -600039852  # N72

$ws.Cells.Item(106, 8).Value = 20692454  # H106
$ws.Cells.Item(106, 9).Value = 26089768  # I106
$ws.Cells.Item(106, 10).Value = 2750  # J106
$ws.Cells.Item(106, 11).Value = 26089768  # K106
$ws.Cells.Item(106, 12).Value = 2750  # L106
$ws.Cells.Item(106, 13).Value = -26089137  # M106
$ws.Cells.Item(106, 14).Value = -4012  # N106

$ws.Cells.Item(138, 8).Value = 2646.25  # H138
$ws.Cells.Item(138, 9).Value = 1370.7778  # I138
$ws.Cells.Item(138, 10).Value = 2926.2317  # J138
$ws.Cells.Item(138, 11).Value = 4112.3334  # K138
$ws.Cells.Item(138, 12).Value = 8778.695099999999  # L138
$ws.Cells.Item(138, 13).Value = 1027.6666  # M138
$ws.Cells.Item(138, 14).Value = -19058.6951  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10294.08  # H32
$ws.Cells.Item(32, 9).Value = 8702.897000000001  # I32
$ws.Cells.Item(32, 11).Value = 8702.897000000001  # K32
$ws.Cells.Item(32, 13).Value = -8415.897000000001  # M32

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 1851.5555  # H5
$ws.Cells.Item(5, 9).Value = 193.33333  # I5
$ws.Cells.Item(5, 10).Value = 3509.7778  # J5
$ws.Cells.Item(5, 11).Value = 193.33333  # K5
$ws.Cells.Item(5, 12).Value = 3509.7778  # L5
$ws.Cells.Item(5, 13).Value = -81.33332999999999  # M5
$ws.Cells.Item(5, 14).Value = -3733.7778  # N5

$ws.Cells.Item(25, 8).Value = 4645.909  # H25
$ws.Cells.Item(25, 9).Value = 1011  # I25
$ws.Cells.Item(25, 10).Value = 9007.799999999999  # J25
$ws.Cells.Item(25, 11).Value = 1011  # K25
$ws.Cells.Item(25, 12).Value = 9007.799999999999  # L25
$ws.Cells.Item(25, 13).Value = -837  # M25
$ws.Cells.Item(25, 14).Value = -9355.799999999999  # N25

$ws.Cells.Item(31, 8).Value = 5344.698  # H31
$ws.Cells.Item(31, 9).Value = 1859.45  # I31
$ws.Cells.Item(31, 10).Value = 6965.744  # J31
$ws.Cells.Item(31, 11).Value = 1859.45  # K31
$ws.Cells.Item(31, 12).Value = 6965.744  # L31
$ws.Cells.Item(31, 13).Value = -1564.45  # M31
$ws.Cells.Item(31, 14).Value = -7555.744  # N31

$ws.Cells.Item(34, 8).Value = 5344.698  # H34
$ws.Cells.Item(34, 9).Value = 1859.45  # I34
$ws.Cells.Item(34, 10).Value = 6965.744  # J34
$ws.Cells.Item(34, 11).Value = 1859.45  # K34
$ws.Cells.Item(34, 12).Value = 6965.744  # L34
$ws.Cells.Item(34, 13).Value = -1657.45  # M34
$ws.Cells.Item(34, 14).Value = -7369.744  # N34

$ws.Cells.Item(75, 8).Value = 70000  # H75
$ws.Cells.Item(75, 10).Value = 70000  # J75
$ws.Cells.Item(75, 12).Value = 70000  # L75
$ws.Cells.Item(75, 14).Value = -71996  # N75

$ws.Cells.Item(78, 8).Value = 70000  # H78
$ws.Cells.Item(78, 10).Value = 70000  # J78
$ws.Cells.Item(78, 12).Value = 210000  # L78
$ws.Cells.Item(78, 14).Value = -219984  # N78

$ws.Cells.Item(109, 8).Value = 30242.5  # H109
$ws.Cells.Item(109, 10).Value = 30242.5  # J109
$ws.Cells.Item(109, 12).Value = 30242.5  # L109
$ws.Cells.Item(109, 14).Value = -32322.5  # N109

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 655.7692  # H5
$ws.Cells.Item(5, 9).Value = 465.9535  # I5
$ws.Cells.Item(5, 10).Value = 1026.7727  # J5
$ws.Cells.Item(5, 11).Value = 1397.8605  # K5
$ws.Cells.Item(5, 12).Value = 3080.3181  # L5
$ws.Cells.Item(5, 13).Value = -1285.8605  # M5
$ws.Cells.Item(5, 14).Value = -3304.3181  # N5

$ws.Cells.Item(64, 8).Value = 1390938  # H64
$ws.Cells.Item(64, 9).Value = 1113.3334  # I64
$ws.Cells.Item(64, 10).Value = 1854212.9  # J64
$ws.Cells.Item(64, 11).Value = 3340.0002  # K64
$ws.Cells.Item(64, 12).Value = 5562638.699999999  # L64
$ws.Cells.Item(64, 13).Value = -3070.0002  # M64
$ws.Cells.Item(64, 14).Value = -5563178.699999999  # N64

$ws.Cells.Item(67, 8).Value = 1390938  # H67
$ws.Cells.Item(67, 9).Value = 1113.3334  # I67
$ws.Cells.Item(67, 10).Value = 1854212.9  # J67
$ws.Cells.Item(67, 11).Value = 3340.0002  # K67
$ws.Cells.Item(67, 12).Value = 5562638.699999999  # L67
$ws.Cells.Item(67, 13).Value = -2404.0002  # M67
$ws.Cells.Item(67, 14).Value = -5564510.699999999  # N67

$ws.Cells.Item(70, 8).Value = 4736.1665  # H70
$ws.Cells.Item(70, 9).Value = 4736.1665  # I70
$ws.Cells.Item(70, 11).Value = 14208.4995  # K70
$ws.Cells.Item(70, 13).Value = -13893.4995  # M70

$ws.Cells.Item(73, 8).Value = 4736.1665  # H73
$ws.Cells.Item(73, 9).Value = 4736.1665  # I73
$ws.Cells.Item(73, 11).Value = 14208.4995  # K73
$ws.Cells.Item(73, 13).Value = -13116.4995  # M73

$ws.Cells.Item(135, 8).Value = 655.7692  # H135
$ws.Cells.Item(135, 9).Value = 465.9535  # I135
$ws.Cells.Item(135, 10).Value = 1026.7727  # J135
$ws.Cells.Item(135, 11).Value = 4193.5815  # K135
$ws.Cells.Item(135, 12).Value = 9240.954299999999  # L135
$ws.Cells.Item(135, 13).Value = -1658.5815  # M135
$ws.Cells.Item(135, 14).Value = -14310.9543  # N135

$ws.Cells.Item(139, 8).Value = 297749.56  # H139
$ws.Cells.Item(139, 9).Value = 478060  # I139
$ws.Cells.Item(139, 10).Value = 6478.846  # J139
$ws.Cells.Item(139, 11).Value = 1434180  # K139
$ws.Cells.Item(139, 12).Value = 19436.538  # L139
$ws.Cells.Item(139, 13).Value = -1429040  # M139
$ws.Cells.Item(139, 14).Value = -29716.538  # N139

$ws.Cells.Item(141, 8).Value = 8925.8125  # H141
$ws.Cells.Item(141, 9).Value = 5164  # I141
$ws.Cells.Item(141, 10).Value = 10635.728  # J141
$ws.Cells.Item(141, 11).Value = 15492  # K141
$ws.Cells.Item(141, 12).Value = 31907.184  # L141
$ws.Cells.Item(141, 13).Value = -10312  # M141
$ws.Cells.Item(141, 14).Value = -42267.18399999999  # N141

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 4000  # H31
$ws.Cells.Item(31, 9).Value = 4000  # I31
$ws.Cells.Item(31, 10).Value = 0  # J31
$ws.Cells.Item(31, 11).Value = 4000  # K31
$ws.Cells.Item(31, 12).Value = 0  # L31
$ws.Cells.Item(31, 13).Value = -3708  # M31
$ws.Cells.Item(31, 14).ClearContents()  # N31

$ws.Cells.Item(37, 8).Value = 4000  # H37
$ws.Cells.Item(37, 9).Value = 4000  # I37
$ws.Cells.Item(37, 10).Value = 0  # J37
$ws.Cells.Item(37, 11).Value = 4000  # K37
$ws.Cells.Item(37, 12).Value = 0  # L37
$ws.Cells.Item(37, 13).Value = -3723  # M37
$ws.Cells.Item(37, 14).ClearContents()  # N37

$ws.Cells.Item(70, 8).Value = 8298.612999999999  # H70
$ws.Cells.Item(70, 9).Value = 9199.478999999999  # I70
$ws.Cells.Item(70, 11).Value = 9199.478999999999  # K70
$ws.Cells.Item(70, 13).Value = -8929.478999999999  # M70

$ws.Cells.Item(73, 8).Value = 8298.612999999999  # H73
$ws.Cells.Item(73, 9).Value = 9199.478999999999  # I73
$ws.Cells.Item(73, 11).Value = 9199.478999999999  # K73
$ws.Cells.Item(73, 13).Value = -8263.478999999999  # M73

$ws.Cells.Item(80, 8).Value = 1568750  # H80
$ws.Cells.Item(80, 9).Value = 3002166.8  # I80
$ws.Cells.Item(80, 11).Value = 3002166.8  # K80
$ws.Cells.Item(80, 13).Value = -3001168.8  # M80

$ws.Cells.Item(83, 8).Value = 1568750  # H83
$ws.Cells.Item(83, 9).Value = 3002166.8  # I83
$ws.Cells.Item(83, 11).Value = 15010834  # K83
$ws.Cells.Item(83, 13).Value = -15005842  # M83

$ws.Cells.Item(132, 8).Value = 30305644  # H132
$ws.Cells.Item(132, 9).Value = 43480270  # I132
$ws.Cells.Item(132, 10).Value = 4002.4  # J132
$ws.Cells.Item(132, 11).Value = 130440810  # K132
$ws.Cells.Item(132, 12).Value = 12007.2  # L132
$ws.Cells.Item(132, 13).Value = -130438280  # M132
$ws.Cells.Item(132, 14).Value = -17067.2  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2734  # H68
$ws.Cells.Item(68, 9).Value = 2002  # I68
$ws.Cells.Item(68, 10).Value = 3100  # J68
$ws.Cells.Item(68, 11).Value = 2002  # K68
$ws.Cells.Item(68, 12).Value = 3100  # L68
$ws.Cells.Item(68, 13).Value = -1253  # M68
$ws.Cells.Item(68, 14).Value = -4598  # N68

$ws.Cells.Item(71, 8).Value = 2734  # H71
$ws.Cells.Item(71, 9).Value = 2002  # I71
$ws.Cells.Item(71, 10).Value = 3100  # J71
$ws.Cells.Item(71, 11).Value = 10010  # K71
$ws.Cells.Item(71, 12).Value = 15500  # L71
$ws.Cells.Item(71, 13).Value = -6266  # M71
$ws.Cells.Item(71, 14).Value = -22988  # N71

$ws.Cells.Item(100, 8).Value = 49826.79  # H100
$ws.Cells.Item(100, 9).Value = 58561.875  # I100
$ws.Cells.Item(100, 10).Value = 3239.6667  # J100
$ws.Cells.Item(100, 11).Value = 58561.875  # K100
$ws.Cells.Item(100, 12).Value = 3239.6667  # L100
$ws.Cells.Item(100, 13).Value = -58020.875  # M100
$ws.Cells.Item(100, 14).Value = -4321.6667  # N100

$ws.Cells.Item(132, 8).Value = 3241.6296  # H132
$ws.Cells.Item(132, 9).Value = 2647  # I132
$ws.Cells.Item(132, 11).Value = 7941  # K132
$ws.Cells.Item(132, 13).Value = -5411  # M132

$ws.Cells.Item(140, 8).Value = 52843.582  # H140
$ws.Cells.Item(140, 10).Value = 52843.582  # J140
$ws.Cells.Item(140, 12).Value = 52843.582  # L140
$ws.Cells.Item(140, 14).Value = -63203.582  # N140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 104720.695  # H62
$ws.Cells.Item(62, 9).Value = 204279.8  # I62
$ws.Cells.Item(62, 10).Value = 42496.25  # J62
$ws.Cells.Item(62, 11).Value = 204279.8  # K62
$ws.Cells.Item(62, 12).Value = 42496.25  # L62
$ws.Cells.Item(62, 13).Value = -203655.8  # M62
$ws.Cells.Item(62, 14).Value = -43744.25  # N62

$ws.Cells.Item(65, 8).Value = 104720.695  # H65
$ws.Cells.Item(65, 9).Value = 204279.8  # I65
$ws.Cells.Item(65, 10).Value = 42496.25  # J65
$ws.Cells.Item(65, 11).Value = 1021399  # K65
$ws.Cells.Item(65, 12).Value = 212481.25  # L65
$ws.Cells.Item(65, 13).Value = -1018279  # M65
$ws.Cells.Item(65, 14).Value = -218721.25  # N65

$ws.Cells.Item(75, 8).Value = 85310  # H75
$ws.Cells.Item(75, 10).Value = 85310  # J75
$ws.Cells.Item(75, 12).Value = 85310  # L75
$ws.Cells.Item(75, 14).Value = -87182  # N75

$ws.Cells.Item(78, 8).Value = 85310  # H78
$ws.Cells.Item(78, 10).Value = 85310  # J78
$ws.Cells.Item(78, 12).Value = 255930  # L78
$ws.Cells.Item(78, 14).Value = -265290  # N78
